$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPUKJVQIJ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T200042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEIaCXVzLXdlc3QtMiJHMEUCIQDALng0WFak1QmBkai7612VQgHa4CozHn7SZvrxohJAIAIgGKQh0dm8SqX9XSdHd13zwEMyvO3%2BAF7foE9MFfFUBhEq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDKoJSfigixRsMYtYyiq6A15Saroy5V6djznxia%2F5xabHOz2OtdhERpBWEq30Ng2WpQOpOgIt5pnGjNDU9a9Avg1pAz%2FjMlD%2BqS6lIr0KvU3SbBFu3mv%2BAm1mauTHDX4wCXg9m5bD2CQ9cSTFypzIl0s6DY%2FYV9iFM%2BTGIDHixJ3VFjF2VIDRRSxqiVE361hEVuIvnp0vkSu3hsMsqN5Zl9pvKdxQtFeGBMNoC0JCv%2BRTsaHZktogD2y1k09HC%2BZjW7U152cB9OcitQdStAjqGDzdnwcQph%2FlojTJqBM8tNR%2FUdI1rCnUvOIH%2Br33fA43Ta4XLiam%2FMm26bhcJZw6AUdJZS07%2BB1yPKLBVKeK1TCWmOt0FJUsKl5aM%2FqSGFyQr68RrZypfgi9wGskWOqzB9tGmVmY7WP9X3Xr7PbN0YKA4l9wqsw6J7cetMoYrY54y7rd9uvgoCQPMOMTGyUVblZRzsplwub46X3qkSRqvHabtKuh9af0xlcdjIt1rnfx6DWdlZTUdhXTH%2FN%2Bz0YH2Ps2W0N4993u1HoZkyqD2oiBGlRwSoqPLl%2B3pzOZfBXPkkZd9W4lu%2B1Buc6UC59%2FFc1WuC%2FfLvkkqQgw7%2Ff8wQY6pQGG5qlC7O0N5dtHOJ72PWbvDZRIqm6s6a2r%2Fpa9ndfcOE8Nre4UuLunoWAOl4%2FSWfL4f%2BoqeC8z%2BNN%2FKQ31%2FD9TR3TMlUrtA6xS0ygEh8HzTwAE%2FSq1JHdmTIM297avH%2Bf1Gyoi2X4mLZtsTwRhAskqpgrQschLtgqPrFAO37ZyASTCobyCCNkxyNg3dHtU5RmR1kbsoIHbLfNfrr2RZVG3vqD29dA%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2004%3A00%3A42%20GMT&X-Amz-Signature=f0fae985171a4a287af55fae6d3aca28eb224ca74597aef1544fded973af1960'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPUKJVQIJ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T200042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEIaCXVzLXdlc3QtMiJHMEUCIQDALng0WFak1QmBkai7612VQgHa4CozHn7SZvrxohJAIAIgGKQh0dm8SqX9XSdHd13zwEMyvO3%2BAF7foE9MFfFUBhEq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDKoJSfigixRsMYtYyiq6A15Saroy5V6djznxia%2F5xabHOz2OtdhERpBWEq30Ng2WpQOpOgIt5pnGjNDU9a9Avg1pAz%2FjMlD%2BqS6lIr0KvU3SbBFu3mv%2BAm1mauTHDX4wCXg9m5bD2CQ9cSTFypzIl0s6DY%2FYV9iFM%2BTGIDHixJ3VFjF2VIDRRSxqiVE361hEVuIvnp0vkSu3hsMsqN5Zl9pvKdxQtFeGBMNoC0JCv%2BRTsaHZktogD2y1k09HC%2BZjW7U152cB9OcitQdStAjqGDzdnwcQph%2FlojTJqBM8tNR%2FUdI1rCnUvOIH%2Br33fA43Ta4XLiam%2FMm26bhcJZw6AUdJZS07%2BB1yPKLBVKeK1TCWmOt0FJUsKl5aM%2FqSGFyQr68RrZypfgi9wGskWOqzB9tGmVmY7WP9X3Xr7PbN0YKA4l9wqsw6J7cetMoYrY54y7rd9uvgoCQPMOMTGyUVblZRzsplwub46X3qkSRqvHabtKuh9af0xlcdjIt1rnfx6DWdlZTUdhXTH%2FN%2Bz0YH2Ps2W0N4993u1HoZkyqD2oiBGlRwSoqPLl%2B3pzOZfBXPkkZd9W4lu%2B1Buc6UC59%2FFc1WuC%2FfLvkkqQgw7%2Ff8wQY6pQGG5qlC7O0N5dtHOJ72PWbvDZRIqm6s6a2r%2Fpa9ndfcOE8Nre4UuLunoWAOl4%2FSWfL4f%2BoqeC8z%2BNN%2FKQ31%2FD9TR3TMlUrtA6xS0ygEh8HzTwAE%2FSq1JHdmTIM297avH%2Bf1Gyoi2X4mLZtsTwRhAskqpgrQschLtgqPrFAO37ZyASTCobyCCNkxyNg3dHtU5RmR1kbsoIHbLfNfrr2RZVG3vqD29dA%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2004%3A00%3A42%20GMT&X-Amz-Signature=242a957d2fd4d37d45bbb0f54c3cba811413cc31145e92fdea79b9ee3281c236'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPUKJVQIJ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T200042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEIaCXVzLXdlc3QtMiJHMEUCIQDALng0WFak1QmBkai7612VQgHa4CozHn7SZvrxohJAIAIgGKQh0dm8SqX9XSdHd13zwEMyvO3%2BAF7foE9MFfFUBhEq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDKoJSfigixRsMYtYyiq6A15Saroy5V6djznxia%2F5xabHOz2OtdhERpBWEq30Ng2WpQOpOgIt5pnGjNDU9a9Avg1pAz%2FjMlD%2BqS6lIr0KvU3SbBFu3mv%2BAm1mauTHDX4wCXg9m5bD2CQ9cSTFypzIl0s6DY%2FYV9iFM%2BTGIDHixJ3VFjF2VIDRRSxqiVE361hEVuIvnp0vkSu3hsMsqN5Zl9pvKdxQtFeGBMNoC0JCv%2BRTsaHZktogD2y1k09HC%2BZjW7U152cB9OcitQdStAjqGDzdnwcQph%2FlojTJqBM8tNR%2FUdI1rCnUvOIH%2Br33fA43Ta4XLiam%2FMm26bhcJZw6AUdJZS07%2BB1yPKLBVKeK1TCWmOt0FJUsKl5aM%2FqSGFyQr68RrZypfgi9wGskWOqzB9tGmVmY7WP9X3Xr7PbN0YKA4l9wqsw6J7cetMoYrY54y7rd9uvgoCQPMOMTGyUVblZRzsplwub46X3qkSRqvHabtKuh9af0xlcdjIt1rnfx6DWdlZTUdhXTH%2FN%2Bz0YH2Ps2W0N4993u1HoZkyqD2oiBGlRwSoqPLl%2B3pzOZfBXPkkZd9W4lu%2B1Buc6UC59%2FFc1WuC%2FfLvkkqQgw7%2Ff8wQY6pQGG5qlC7O0N5dtHOJ72PWbvDZRIqm6s6a2r%2Fpa9ndfcOE8Nre4UuLunoWAOl4%2FSWfL4f%2BoqeC8z%2BNN%2FKQ31%2FD9TR3TMlUrtA6xS0ygEh8HzTwAE%2FSq1JHdmTIM297avH%2Bf1Gyoi2X4mLZtsTwRhAskqpgrQschLtgqPrFAO37ZyASTCobyCCNkxyNg3dHtU5RmR1kbsoIHbLfNfrr2RZVG3vqD29dA%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2004%3A00%3A42%20GMT&X-Amz-Signature=574a044b34d7c4e35ebc47313fdb4794361b9821ff77f5083f0ba024f34f0d8e'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPUKJVQIJ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T200042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEIaCXVzLXdlc3QtMiJHMEUCIQDALng0WFak1QmBkai7612VQgHa4CozHn7SZvrxohJAIAIgGKQh0dm8SqX9XSdHd13zwEMyvO3%2BAF7foE9MFfFUBhEq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDKoJSfigixRsMYtYyiq6A15Saroy5V6djznxia%2F5xabHOz2OtdhERpBWEq30Ng2WpQOpOgIt5pnGjNDU9a9Avg1pAz%2FjMlD%2BqS6lIr0KvU3SbBFu3mv%2BAm1mauTHDX4wCXg9m5bD2CQ9cSTFypzIl0s6DY%2FYV9iFM%2BTGIDHixJ3VFjF2VIDRRSxqiVE361hEVuIvnp0vkSu3hsMsqN5Zl9pvKdxQtFeGBMNoC0JCv%2BRTsaHZktogD2y1k09HC%2BZjW7U152cB9OcitQdStAjqGDzdnwcQph%2FlojTJqBM8tNR%2FUdI1rCnUvOIH%2Br33fA43Ta4XLiam%2FMm26bhcJZw6AUdJZS07%2BB1yPKLBVKeK1TCWmOt0FJUsKl5aM%2FqSGFyQr68RrZypfgi9wGskWOqzB9tGmVmY7WP9X3Xr7PbN0YKA4l9wqsw6J7cetMoYrY54y7rd9uvgoCQPMOMTGyUVblZRzsplwub46X3qkSRqvHabtKuh9af0xlcdjIt1rnfx6DWdlZTUdhXTH%2FN%2Bz0YH2Ps2W0N4993u1HoZkyqD2oiBGlRwSoqPLl%2B3pzOZfBXPkkZd9W4lu%2B1Buc6UC59%2FFc1WuC%2FfLvkkqQgw7%2Ff8wQY6pQGG5qlC7O0N5dtHOJ72PWbvDZRIqm6s6a2r%2Fpa9ndfcOE8Nre4UuLunoWAOl4%2FSWfL4f%2BoqeC8z%2BNN%2FKQ31%2FD9TR3TMlUrtA6xS0ygEh8HzTwAE%2FSq1JHdmTIM297avH%2Bf1Gyoi2X4mLZtsTwRhAskqpgrQschLtgqPrFAO37ZyASTCobyCCNkxyNg3dHtU5RmR1kbsoIHbLfNfrr2RZVG3vqD29dA%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2004%3A00%3A42%20GMT&X-Amz-Signature=eb3f71db89f9ac06ded4804ab7af2c1ad1376687e228c4a3a28c2d34001629a2'
